$wb = $excel.ActiveWorkbook

# Sheet 1: Narrator Votes Averages
$ws1 = $wb.Worksheets.Item("Narrator Votes Averages")
$ws1.Range("B2").Value = 15
$ws1.Range("C2").Value = 27.5
$ws1.Range("B3").Value = 15.66666666666667
$ws1.Range("C3").Value = 30.27777777777777
$ws1.Range("B4").Value = 69.33333333333333
$ws1.Range("C4").Value = 42.22222222222221

# Sheet 2: Votes Not Narrator Averages
$ws2 = $wb.Worksheets.Item("Votes Not Narrator Averages")
$ws2.Range("B2").Value = 38.45454545454545
$ws2.Range("C2").Value = 49.07407407407407

# Sheet 3: Correct Votes Averages
$ws3 = $wb.Worksheets.Item("Correct Votes Averages")
$ws3.Range("B2").Value = 42.90909090909091
$ws3.Range("C2").Value = 50.98605098605099
